$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-51: Coin (B), Link (C), Price (D), Volume1h (E)
$data = @(
    @{B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='30.088.33'; E='  -1.83%  '},
    @{B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='1.829.04'; E='  -3.33%  '},
    @{B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='1.000'; E='  -0.06%  '},
    @{B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='229.60'; E='  -3.77%  '},
    @{B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='1.000'; E='  -0.03%  '},
    @{B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='0.4621'; E='  -4.42%  '},
    @{B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.2696'; E='  -5.87%  '},
    @{B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.06227'; E='  -5.00%  '},
    @{B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.842.85'; E='  -4.56%  '},
    @{B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.07362'; E='  -1.49%  '},
    @{B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='16.03'; E='  -4.17%  '},
    @{B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='4.890'; E='  -4.25%  '},
    @{B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='82.87'; E='  -5.97%  '},
    @{B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='0.6179'; E='  -7.40%  '},
    @{B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='30.025.34'; E='  -1.95%  '},
    @{B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.000'; E='  -0.02%  '},
    @{B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='227.04'; E='  -1.46%  '},
    @{B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000007255'; E='  -4.26%  '},
    @{B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.001'; E='  +0.01%  '},
    @{B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='12.30'; E='  -7.35%  '},
    @{B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='4.822'; E='  -8.57%  '},
    @{B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='5.803'; E='  -6.48%  '},
    @{B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='165.30'; E='  -2.69%  '},
    @{B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='9.132'; E='  -2.48%  '},
    @{B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='17.73'; E='  -5.90%  '},
    @{B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='1.835'; E='  -6.46%  '},
    @{B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.1008'; E='  -1.62%  '},
    @{B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='1.367'; E='  -2.35%  '},
    @{B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='4.028'; E='  -7.07%  '},
    @{B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='3.741'; E='  -7.09%  '},
    @{B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.04778'; E='  -5.55%  '},
    @{B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.121'; E='  -7.74%  '},
    @{B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.6968'; E='  -7.69%  '},
    @{B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='2.686'; E='  -0.96%  '},
    @{B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.01801'; E='  -3.81%  '},
    @{B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.599'; E='  -1.76%  '},
    @{B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='0.8881'; E='  -3.61%  '},
    @{B='PaxDollar'; C='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D='0.9999'; E='  -0.30%  '},
    @{B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='1.906'; E='  -7.73%  '},
    @{B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='102.51'; E='  -4.22%  '},
    @{B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='5.466'; E='  -3.57%  '},
    @{B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.3972'; E='  -7.41%  '},
    @{B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='6.892'; E='  -7.12%  '},
    @{B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1185'; E='  -7.10%  '},
    @{B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='59.09'; E='  -8.37%  '},
    @{B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='8.407'; E='  -6.45%  '},
    @{B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.05516'; E='  -2.58%  '},
    @{B='Elrond'; C='https://coinranking.com/coin/omwkOTglq+elrond-egld'; D='32.44'; E='  -4.37%  '},
    @{B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='1.345'; E='  -9.98%  '},
    @{B='Decentraland'; C='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D='0.3613'; E='  -7.18%  '}
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 2).Value2 = $item.B
    $ws.Cells.Item($row, 3).Value2 = $item.C

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value2 = $item.D
    $dCell.Style = "Normal"

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value2 = $item.E
    $eCell.Style = "Normal"

    $row = $row + 1
}
